# Append 24 new transaction rows (rows 80-103) for week ending 2021-01-31
# to the "Konto" sheet, mirroring the existing Datum/Receipt Number/Konto/
# Beskrivning/Debet/Kredit column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the date NumberFormat already applied to column A's existing rows.
$dateFormat = $ws.Cells.Item(2, 1).NumberFormat

# Each entry: Datum (serial), Receipt Number, Konto, Beskrivning, Debet, Kredit
# A $null value for Receipt Number / Debet / Kredit means "leave blank".
$newRows = @(
    @(44221, 1252204, 3011, "Order 1252204 Swish +46730393329", $null, 502.68),
    @(44221, 1252204, 2611, "Order 1252204 Swish +46730393329", $null, 60.32),
    @(44221, 1252204, 1930, "Order 1252204 Swish +46730393329", 563, $null),

    @(44222, $null, 1220, "ELGIGANTEN STHL K0135", 808.8, $null),
    @(44222, $null, 2641, "ELGIGANTEN STHL K0135", 202.2, $null),
    @(44222, $null, 1930, "ELGIGANTEN STHL K0135", $null, 1011),

    @(44222, 8261903, 3011, "Order 8261903 Swish +46733304498", $null, 691.0700000000001),
    @(44222, 8261903, 2611, "Order 8261903 Swish +46733304498", $null, 82.93000000000001),
    @(44222, 8261903, 1930, "Order 8261903 Swish +46733304498", 774, $null),

    @(44227, 4311427, 3011, "Order 4311427 Swish +46727242898", $null, 806.25),
    @(44227, 4311427, 2611, "Order 4311427 Swish +46727242898", $null, 96.75),
    @(44227, 4311427, 1930, "Order 4311427 Swish +46727242898", 903, $null),

    @(44227, 5311815, 3011, "Order 5311815 Swish +46708751433", $null, 616.0700000000001),
    @(44227, 5311815, 2611, "Order 5311815 Swish +46708751433", $null, 73.93000000000001),
    @(44227, 5311815, 1930, "Order 5311815 Swish +46708751433", 690, $null),

    @(44227, 2311908, 3011, "Order 2311908 Swish +46709526084", $null, 739.29),
    @(44227, 2311908, 2611, "Order 2311908 Swish +46709526084", $null, 88.70999999999999),
    @(44227, 2311908, 1930, "Order 2311908 Swish +46709526084", 828, $null),

    @(44227, 9311931, 3011, "Order 9311931 Swish +46706649892", $null, 691.0700000000001),
    @(44227, 9311931, 2611, "Order 9311931 Swish +46706649892", $null, 82.93000000000001),
    @(44227, 9311931, 1930, "Order 9311931 Swish +46706649892", 774, $null),

    @(44227, 3312155, 3011, "Order 3312155 Swish +46707676358", $null, 1106.25),
    @(44227, 3312155, 2611, "Order 3312155 Swish +46707676358", $null, 132.75),
    @(44227, 3312155, 1930, "Order 3312155 Swish +46707676358", 1239, $null)
)

$startRow = 80
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $datum = $data[0]
    $receiptNumber = $data[1]
    $konto = $data[2]
    $beskrivning = $data[3]
    $debet = $data[4]
    $kredit = $data[5]

    $ws.Cells.Item($r, 1).Value = $datum
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat

    if ($null -eq $receiptNumber) {
        $ws.Cells.Item($r, 2).Value = ""
    } else {
        $ws.Cells.Item($r, 2).Value = $receiptNumber
    }

    $ws.Cells.Item($r, 3).Value = $konto
    $ws.Cells.Item($r, 4).Value = $beskrivning

    if ($null -eq $debet) {
        $ws.Cells.Item($r, 5).Value = ""
    } else {
        $ws.Cells.Item($r, 5).Value = $debet
    }

    if ($null -eq $kredit) {
        $ws.Cells.Item($r, 6).Value = ""
    } else {
        $ws.Cells.Item($r, 6).Value = $kredit
    }
}
